$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "22.346.49"
Set-TextValue $ws.Range("E2") "  -0.40%  "

Set-TextValue $ws.Range("D3") "1.564.64"
Set-TextValue $ws.Range("E3") "  -0.59%  "

Set-TextValue $ws.Range("D4") "1.004"
Set-TextValue $ws.Range("E4") "  +0.35%  "

Set-TextValue $ws.Range("E5") "  +0.27%  "

Set-TextValue $ws.Range("D6") "290.59"
Set-TextValue $ws.Range("E6") "  +0.09%  "

Set-TextValue $ws.Range("D7") "0.3760"
Set-TextValue $ws.Range("E7") "  +2.05%  "

Set-TextValue $ws.Range("D8") "49.18"
Set-TextValue $ws.Range("E8") "  -0.27%  "

Set-TextValue $ws.Range("D9") "0.3395"
Set-TextValue $ws.Range("E9") "  -0.38%  "

Set-TextValue $ws.Range("D10") "0.07566"
Set-TextValue $ws.Range("E10") "  -1.46%  "

Set-TextValue $ws.Range("D11") "1.133"
Set-TextValue $ws.Range("E11") "  -3.56%  "

Set-TextValue $ws.Range("D12") "1.004"
Set-TextValue $ws.Range("E12") "  +0.36%  "

Set-TextValue $ws.Range("D13") "20.97"
Set-TextValue $ws.Range("E13") "  -1.90%  "

Set-TextValue $ws.Range("D14") "5.951"
Set-TextValue $ws.Range("E14") "  -2.07%  "

Set-TextValue $ws.Range("D15") "6.896"
Set-TextValue $ws.Range("E15") "  -0.72%  "

Set-TextValue $ws.Range("D16") "1.565.43"
Set-TextValue $ws.Range("E16") "  -0.02%  "

Set-TextValue $ws.Range("D17") "0.00001128"
Set-TextValue $ws.Range("E17") "  -0.54%  "

Set-TextValue $ws.Range("D18") "89.68"
Set-TextValue $ws.Range("E18") "  -0.80%  "

Set-TextValue $ws.Range("E19") "  -0.14%  "

Set-TextValue $ws.Range("D20") "1.003"
Set-TextValue $ws.Range("E20") "  +0.27%  "

Set-TextValue $ws.Range("D21") "16.54"
Set-TextValue $ws.Range("E21") "  -0.34%  "

Set-TextValue $ws.Range("D22") "6.184"
Set-TextValue $ws.Range("E22") "  -1.54%  "

Set-TextValue $ws.Range("D23") "11.92"
Set-TextValue $ws.Range("E23") "  -1.14%  "

Set-TextValue $ws.Range("D24") "22.340.91"
Set-TextValue $ws.Range("E24") "  -0.46%  "

Set-TextValue $ws.Range("D25") "2.378"
Set-TextValue $ws.Range("E25") "  +0.79%  "

Set-TextValue $ws.Range("D26") "2.694"
Set-TextValue $ws.Range("E26") "  -6.39%  "

Set-TextValue $ws.Range("D27") "20.13"
Set-TextValue $ws.Range("E27") "  -0.04%  "

Set-TextValue $ws.Range("D28") "146.97"
Set-TextValue $ws.Range("E28") "  +0.10%  "

Set-TextValue $ws.Range("D29") "5.016"
Set-TextValue $ws.Range("E29") "  +0.63%  "

Set-TextValue $ws.Range("D30") "125.64"
Set-TextValue $ws.Range("E30") "  -0.14%  "

Set-TextValue $ws.Range("D31") "1.741.38"
Set-TextValue $ws.Range("E31") "  +0.00%  "

Set-TextValue $ws.Range("D32") "2.017"
Set-TextValue $ws.Range("E32") "  -0.07%  "

Set-TextValue $ws.Range("D33") "0.9853"
Set-TextValue $ws.Range("E33") "  -3.84%  "

Set-TextValue $ws.Range("D34") "6.018"
Set-TextValue $ws.Range("E34") "  -3.69%  "

Set-TextValue $ws.Range("D35") "10.04"
Set-TextValue $ws.Range("E35") "  -0.51%  "

Set-TextValue $ws.Range("D36") "1.426"
Set-TextValue $ws.Range("E36") "  +9.73%  "

Set-TextValue $ws.Range("D37") "0.08478"
Set-TextValue $ws.Range("E37") "  -0.03%  "

Set-TextValue $ws.Range("D38") "0.02483"
Set-TextValue $ws.Range("E38") "  -2.75%  "

Set-TextValue $ws.Range("D39") "0.2288"
Set-TextValue $ws.Range("E39") "  -1.64%  "

Set-TextValue $ws.Range("D40") "0.06439"
Set-TextValue $ws.Range("E40") "  -0.76%  "

Set-TextValue $ws.Range("D41") "5.389"
Set-TextValue $ws.Range("E41") "  -3.12%  "

Set-TextValue $ws.Range("D42") "0.6280"
Set-TextValue $ws.Range("E42") "  -1.43%  "

Set-TextValue $ws.Range("D43") "11.23"
Set-TextValue $ws.Range("E43") "  -4.56%  "

Set-TextValue $ws.Range("E44") "  +0.24%  "

Set-TextValue $ws.Range("D45") "13.86"
Set-TextValue $ws.Range("E45") "  -3.10%  "

Set-TextValue $ws.Range("D46") "3.798"
Set-TextValue $ws.Range("E46") "  +0.93%  "

Set-TextValue $ws.Range("D47") "0.5903"
Set-TextValue $ws.Range("E47") "  -1.72%  "

Set-TextValue $ws.Range("D48") "2.064"
Set-TextValue $ws.Range("E48") "  -2.65%  "

Set-TextValue $ws.Range("E49") "  -0.71%  "

Set-TextValue $ws.Range("D50") "124.25"
Set-TextValue $ws.Range("E50") "  -0.49%  "

Set-TextValue $ws.Range("D51") "0.07318"
Set-TextValue $ws.Range("E51") "  +0.30%  "

